# issue #5: stock data from json to db
# Adds three new columns (category, source_file, index) to the 股票 (stock)
# sheet, shifting the existing property_category/date/legislator_name/
# legislator_id block right and populating the new columns with data that
# mirrors the source-file/category metadata now tracked alongside each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)  # 股票 (stock) sheet

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data row 2 -------------------------------------------------------------
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-04-12"
$ws.Range("K2").Value = "王育敏"
$ws.Range("L2").Value = 1728
$ws.Range("M2").Value = "tmp48bc1"
$ws.Range("N2").Value = 51

# --- Data row 3 -------------------------------------------------------------
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2012-04-12"
$ws.Range("K3").Value = "王育敏"
$ws.Range("L3").Value = 1728
$ws.Range("M3").Value = "tmp48bc1"
$ws.Range("N3").Value = 52
